$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove header row styling (bold font, thin border, centered/top alignment)
#    so header cells revert to the default (unstyled) cell format.
$ws.Range("A1:AA1").ClearFormats()

# 2. Clear the old "Unnamed: 0" label from A1 (header for the index column).
$ws.Range("A1").ClearContents()

# 3. Drop the "arg3" (D) and "literal" (M) columns' data rows 3-8 -- the
#    data-cleaning fix excludes these columns from the aggregated stats.
$ws.Range("D3:D8").ClearContents()
$ws.Range("M3:M8").ClearContents()

# 4. Updated aggregate statistics (recomputed after excluding arg3/literal
#    columns from the pre/post/total fixation data).
$ws.Range("B3").Value = 0
$ws.Range("F3").Value = 35
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = 5
$ws.Range("K3").Value = 12
$ws.Range("O3").Value = 28
$ws.Range("P3").Value = 17
$ws.Range("Q3").Value = 6
$ws.Range("T3").Value = 17
$ws.Range("V3").Value = 33
$ws.Range("X3").Value = 12
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 105
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 24
$ws.Range("I4").Value = 8
$ws.Range("K4").Value = 20
$ws.Range("O4").Value = 59
$ws.Range("P4").Value = 29
$ws.Range("Q4").Value = 8
$ws.Range("T4").Value = 38
$ws.Range("V4").Value = 136
$ws.Range("X4").Value = 15
$ws.Range("B5").Value = 133.35
$ws.Range("F5").Value = 42930.67
$ws.Range("G5").Value = 9893.76
$ws.Range("H5").Value = 7208.49
$ws.Range("I5").Value = 2511.12
$ws.Range("K5").Value = 11495.4
$ws.Range("O5").Value = 23125.58
$ws.Range("P5").Value = 18294.78
$ws.Range("Q5").Value = 3854.22
$ws.Range("T5").Value = 20321.23
$ws.Range("V5").Value = 94349.39
$ws.Range("X5").Value = 7215.56
$ws.Range("B6").Value = 0.06
$ws.Range("C6").Value = 1.25
$ws.Range("E6").Value = 4.03
$ws.Range("F6").Value = 20.11
$ws.Range("G6").Value = 4.63
$ws.Range("H6").Value = 3.38
$ws.Range("I6").Value = 1.18
$ws.Range("J6").Value = 3.94
$ws.Range("K6").Value = 5.38
$ws.Range("N6").Value = 1.02
$ws.Range("O6").Value = 10.83
$ws.Range("P6").Value = 8.57
$ws.Range("Q6").Value = 1.81
$ws.Range("R6").Value = 0.65
$ws.Range("T6").Value = 9.52
$ws.Range("V6").Value = 44.19
$ws.Range("W6").Value = 2.58
$ws.Range("X6").Value = 3.38
$ws.Range("Y6").Value = 1.7
$ws.Range("AA6").Value = 0.12
$ws.Range("B7").Value = 133.35
$ws.Range("F7").Value = 408.86
$ws.Range("G7").Value = 618.36
$ws.Range("H7").Value = 300.35
$ws.Range("I7").Value = 313.89
$ws.Range("K7").Value = 574.77
$ws.Range("O7").Value = 391.96
$ws.Range("P7").Value = 630.85
$ws.Range("Q7").Value = 481.78
$ws.Range("T7").Value = 534.77
$ws.Range("V7").Value = 693.75
$ws.Range("X7").Value = 481.04

# 5. Remove the two trailing blank rows (10 and 11) left over from the
#    stale export.
$ws.Rows("10:11").Delete()
